# Edit script: apply "Alterações e Correções sugeridas" changes
#  - Alteração da Visão Geral
#  - Alterações Requisitos e Descrição dos UC
#  - Alterações Especificações UC
#  - Melhorias nos Códigos e Interfaces

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("REQUISITOS DO SISTEMA")
$ws2 = $wb.Worksheets.Item("CASOS DE USO")

# -------------------------------------------------------------------
# Write the brand-new text values first (in the same order the source
# workbook introduces them) so the shared string table grows in a
# natural, stable order, matching how the content was authored.
# -------------------------------------------------------------------
$ws1.Range("D11").Value = "MEDICO/RECEPCIONISTA/DIRETOR"
$ws1.Range("D12").Value = "RECEPCIONISTA/MEDICO/DIRETOR"
$ws1.Range("B4").Value  = "O diretor irá cadastrar o usuário no sistema."
$ws1.Range("B5").Value  = "Enviar e-mail para o usuário contendo link de ativação."
$ws1.Range("B11").Value = "O usuário poderá fazer marcação, desmarcação, listagem de consultas, reagendamento"
$ws1.Range("B12").Value = "O usuário irá manter os pacientes."

# -------------------------------------------------------------------
# Sheet 1: REQUISITOS DO SISTEMA
# -------------------------------------------------------------------

$ws1.Range("A1").Value = "REQUISITOS DO SISTEMA"

# Row 3 - RE01
$ws1.Range("A3").Value = "RE01"
$ws1.Range("B3").Value = "O sistema será web e irá rodar na intranet."
$ws1.Range("C3").Value = "-"
$ws1.Range("D3").Value = "-"

# Row 4 - RE02
$ws1.Range("A4").Value = "RE02"
$ws1.Range("C4").Value = "Registrar Usuário"
$ws1.Range("D4").Value = "DIRETOR"

# Row 5 - RE03
$ws1.Range("A5").Value = "RE03"
$ws1.Range("C5").Value = "Enviar e-mail de ativação"
$ws1.Range("D5").Value = "SISTEMA"

# Row 6 - RE04
$ws1.Range("A6").Value = "RE04"
$ws1.Range("B6").Value = "O usuário deverá clicar no link de ativação, ativando assim o seu status para ativo."
$ws1.Range("C6").Value = "Ativar perfil"
$ws1.Range("D6").Value = "MEDICO/RECEPCIONISTA"

# Row 7 - RE05
$ws1.Range("A7").Value = "RE05"
$ws1.Range("B7").Value = "Fornecer opção para recuperação de senha"
$ws1.Range("C7").Value = "Recuperar Senha"
$ws1.Range("D7").Value = "MEDICO/RECEPCIONISTA/ DIRETOR"

# Row 8 - RE06
$ws1.Range("A8").Value = "RE06"
$ws1.Range("B8").Value = "Permitir que os usuários possam alterar seus dados"
$ws1.Range("C8").Value = "Alterar Dados de Usuário"
$ws1.Range("D8").Value = "MEDICO/RECEPCIONISTA/ DIRETOR"

# Row 9 - RE07
$ws1.Range("A9").Value = "RE07"
$ws1.Range("B9").Value = "O usuário precisa se logar para utilizar o sistema"
$ws1.Range("C9").Value = "Fazer Login"
$ws1.Range("D9").Value = "MEDICO/RECEPCIONISTA/ DIRETOR"

# Row 10 - RE08
$ws1.Range("A10").Value = "RE08"
$ws1.Range("B10").Value = "O diretor poderá listar e excluir usuários"
$ws1.Range("C10").Value = "Manter usuário"
$ws1.Range("D10").Value = "DIRETOR"

# Row 11 - RE09
$ws1.Range("A11").Value = "RE09"
$ws1.Range("C11").Value = "Manter consulta"

# Row 12 - RE10
$ws1.Range("A12").Value = "RE10"
$ws1.Range("C12").Value = "Manter paciente"

# Row 13 - RE11
$ws1.Range("A13").Value = "RE11"
$ws1.Range("B13").Value = "O médico poderá realizar marcação de ""consultas de encaixe"" (consulta sem agendamento)"
$ws1.Range("C13").Value = "Manter consulta de encaixe"
$ws1.Range("D13").Value = "MEDICO/DIRETOR"

# Row 14 - RE12
$ws1.Range("A14").Value = "RE12"
$ws1.Range("B14").Value = "O médico poderá criar e alterar prontuário de um paciente."
$ws1.Range("C14").Value = "Manter prontuário"
$ws1.Range("D14").Value = "MEDICO/DIRETOR"

# Row 15 - RE13
$ws1.Range("A15").Value = "RE13"
$ws1.Range("B15").Value = "O médico poderá solicitar exames."
$ws1.Range("C15").Value = "Solicitar Exame"
$ws1.Range("D15").Value = "MEDICO/DIRETOR"

# Column D width: 43 -> 48.140625 (closest achievable granularity in this engine)
$ws1.Columns.Item(4).ColumnWidth = 47.26

# -------------------------------------------------------------------
# Sheet 2: CASOS DE USO
# -------------------------------------------------------------------

$ws2.Range("A2").Value = "UC"

# Row 3 - UC01
$ws2.Range("A3").Value = "UC01"
$ws2.Range("B3").Value = "Registrar Usuário"
$ws2.Range("C3").Value = "DIRETOR"

# Row 4 - UC02
$ws2.Range("A4").Value = "UC02"
$ws2.Range("B4").Value = "Enviar e-mail de ativação"
$ws2.Range("C4").Value = "SISTEMA"

# Row 5 - UC03
$ws2.Range("A5").Value = "UC03"
$ws2.Range("B5").Value = "Ativar usuário"
$ws2.Range("C5").Value = "MEDICO/RECEPCIONISTA"

# Row 6 - UC04
$ws2.Range("A6").Value = "UC04"
$ws2.Range("B6").Value = "Recuperar Senha"
$ws2.Range("C6").Value = "MEDICO/RECEPCIONISTA/ DIRETOR"

# Row 7 - UC05
$ws2.Range("A7").Value = "UC05"
$ws2.Range("B7").Value = "Fazer Login"
$ws2.Range("C7").Value = "MEDICO/RECEPCIONISTA/ DIRETOR"

# Row 8 - UC06
$ws2.Range("A8").Value = "UC06"
$ws2.Range("B8").Value = "Alterar Dados do Usuario"
$ws2.Range("C8").Value = "MEDICO/RECEPCIONISTA/ DIRETOR"

# Row 9 - UC07
$ws2.Range("A9").Value = "UC07"
$ws2.Range("B9").Value = "Manter usuário"
$ws2.Range("C9").Value = "DIRETOR"

# Row 10 - UC08
$ws2.Range("A10").Value = "UC08"
$ws2.Range("B10").Value = "Manter consulta"
$ws2.Range("C10").Value = "MEDICO/RECEPCIONISTA/DIRETOR"

# Row 11 - UC09
$ws2.Range("A11").Value = "UC09"
$ws2.Range("B11").Value = "Manter paciente"
$ws2.Range("C11").Value = "MEDICO/RECEPCIONISTA/DIRETOR"

# Row 12 - UC10
$ws2.Range("A12").Value = "UC10"
$ws2.Range("B12").Value = "Manter consulta de encaixe"
$ws2.Range("C12").Value = "MEDICO/DIRETOR"

# Row 13 - UC11
$ws2.Range("A13").Value = "UC11"
$ws2.Range("B13").Value = "Manter prontuário"
$ws2.Range("C13").Value = "MEDICO/DIRETOR"

# Row 14 - UC12
$ws2.Range("A14").Value = "UC12"
$ws2.Range("B14").Value = "Solicitar Exame"
$ws2.Range("C14").Value = "MEDICO/DIRETOR"

# Column C width: 42.5703125 -> 50
$ws2.Columns.Item(3).ColumnWidth = 49.2

# -------------------------------------------------------------------
# Selections (set sheet1 first, sheet2 last so sheet2 stays the active tab)
# -------------------------------------------------------------------
$ws1.Range("D13").Select()
$ws2.Range("B10").Select()
